$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("輸出")

for ($i = 1; $i -le 6; $i++) {
    $after = $wb.Worksheets.Item($wb.Worksheets.Count)
    $src.Copy($null, $after)
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = "輸出$i"
}
